$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title "Team " + "Moneyball" -> single run "Team Moneyball"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$r1 = $d.Range($p1.Start, $p1.End - 1)
$r1.Text = "Team  Moneyball"
$r1b = $d.Range($p1.Start, $d.Paragraphs(1).Range.End - 1)
$r1b.Text = "Team Moneyball"

# ---------------------------------------------------------------------------
# 2) "Project Name: " / "...Analysis" / " " -> single run (text unchanged)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$r3 = $d.Range($p3.Start, $p3.End - 1)
$r3.Text = "Project Name: Press Sports App Retention and Engagement Analysis  "
$r3b = $d.Range($p3.Start, $d.Paragraphs(3).Range.End - 1)
$r3b.Text = "Project Name: Press Sports App Retention and Engagement Analysis "

# ---------------------------------------------------------------------------
# 3) "Project Description: ..." multiple runs -> single run (text unchanged)
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4).Range
$r4 = $d.Range($p4.Start, $p4.End - 1)
$r4.Text = "Project Description:  Press Sports App users share an environment to engage with others who share similar interests. We are searching for a method to analyze and understand the driving factors of engagement with the app and other users around the country."
$r4b = $d.Range($p4.Start, $d.Paragraphs(4).Range.End - 1)
$r4b.Text = "Project Description: Press Sports App users share an environment to engage with others who share similar interests. We are searching for a method to analyze and understand the driving factors of engagement with the app and other users around the country."

# ---------------------------------------------------------------------------
# 4) Replace paragraph 5 ("Research questions to be asked...") with the new
#    "Research Question to be Asked" text, then build out the whole
#    hypothesis-testing section as new paragraphs after it.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5).Range
$r5 = $d.Range($p5.Start, $p5.End - 1)
$r5.Text = "Research Question to be Asked: What is the type of user engagement that drives the highest level of user interaction?"

# Paragraph 6 is still the old "Data set to be used: user_data.csv" -- remove
# it now (its replacement is appended at the end of the new block below).
$d.Paragraphs(6).Range.Delete()

$newParas = @(
    @{ Text = "Test to check total number of students at school against number of posts"; Bold = $true },
    @{ Text = "H0 = R^2(total number of students at a school: number of posts) = 0"; Bold = $false },
    @{ Text = "Ha1 = R^2(total number of students at a school: number of posts) > 0"; Bold = $false },
    @{ Text = "Ha2 = R^2(total number of students at a school: number of posts) < 0"; Bold = $false },
    @{ Text = "Test to check total number of students on a club against number of actions"; Bold = $true },
    @{ Text = "H0 = R^2(total number of students on a club: number of actions) = 0"; Bold = $false },
    @{ Text = "Hb1 = R^2(total number of students on a club: number of actions) > 0"; Bold = $false },
    @{ Text = "Hb2 = R^2(total number of students on a club: number of action) < 0"; Bold = $false },
    @{ Text = "Test to check number of users followed against number of posts"; Bold = $true },
    @{ Text = "H0 = R^2(number of users followed on suggest feed: number of posts) = 0"; Bold = $false },
    @{ Text = "Hc1 = R^2(number of users followed on suggest feed: number of posts) > 0"; Bold = $false },
    @{ Text = "Hc2 = R^2(number of users followed on suggest feed: number of posts) < 0"; Bold = $false },
    @{ Text = "Test to check number of users followed against number of actions"; Bold = $true },
    @{ Text = "H0 = R^2(number of users followed on suggest feed: number of actions) = 0"; Bold = $false },
    @{ Text = "Hd1 = R^2(number of users followed on suggest feed: number of actions) > 0"; Bold = $false },
    @{ Text = "Hd2 = R^2(number of users followed on suggest feed: number of actions) < 0"; Bold = $false },
    @{ Text = "Test to check number of post likes against number of posts"; Bold = $true },
    @{ Text = "H0 = R^2(number of post likes: number of posts) = 0"; Bold = $false },
    @{ Text = "He1 = R^2(number of post likes: number of posts) > 0"; Bold = $false },
    @{ Text = "He2 = R^2(number of post likes: number of posts) < 0"; Bold = $false },
    @{ Text = ""; Bold = $true },
    @{ Text = "Test to check number of post likes against number of actions"; Bold = $true },
    @{ Text = "H0 = R^2(number of post likes: number of actions) = 0"; Bold = $false },
    @{ Text = "Hf1 = R^2(number of post likes: number of actions) > 0"; Bold = $false },
    @{ Text = "Hf2 = R^2(number of post likes: number of actions) < 0"; Bold = $false },
    @{ Text = "Data Set to be Used: user_data.csv"; Bold = $false },
    @{ Text = ""; Bold = $false }
)

# First pass: create all of the (still plain/unformatted) empty paragraphs
# that will hold the new content, so that none of them can inherit bold
# paragraph-mark formatting from a previous sibling.
$anchorIdx = 5
for ($k = 0; $k -lt $newParas.Count; $k++) {
    $d.Paragraphs($anchorIdx + $k).Range.InsertParagraphAfter()
}

# Second pass: fill in text + bold, now that every paragraph mark in the
# range is still using the default (non-bold) run formatting.
$idx = 6
foreach ($item in $newParas) {
    if ($item.Text -ne "") {
        $para = $d.Paragraphs($idx).Range
        $para.Text = "~" + $item.Text
        $para2 = $d.Paragraphs($idx).Range
        $para2.Text = $item.Text
    }
    if ($item.Bold) {
        $d.Paragraphs($idx).Range.Font.Bold = $true
        $d.Paragraphs($idx).Range.Font.BoldBi = $true
    }
    $idx = $idx + 1
}

# ---------------------------------------------------------------------------
# 5) "Data Collection" / ": Drew Williams" -> single run (text unchanged)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Data Collection: Drew Williams`r") {
        $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
        $rng.Text = "Data Collection:  Drew Williams"
        $rng2 = $d.Range($p.Range.Start, $p.Range.End - 1)
        $rng2.Text = "Data Collection: Drew Williams"
        break
    }
}

# ---------------------------------------------------------------------------
# 6) "PowerPoint: " / "Team" -> single run (text unchanged), and remove the
#    two trailing empty paragraphs at the very end of the document.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "PowerPoint: Team`r") {
        $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
        $rng.Text = "PowerPoint:  Team"
        $rng2 = $d.Range($p.Range.Start, $p.Range.End - 1)
        $rng2.Text = "PowerPoint: Team"
        break
    }
}

$n = $d.Paragraphs.Count
$pn = $d.Paragraphs($n).Range
$pn1 = $d.Paragraphs($n - 1).Range
$mark = $d.Range($pn1.End - 1, $pn.End - 1)
$mark.Delete()

$n = $d.Paragraphs.Count
$pn = $d.Paragraphs($n).Range
$pn1 = $d.Paragraphs($n - 1).Range
$mark = $d.Range($pn1.End - 1, $pn.End - 1)
$mark.Delete()
